$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.022.60"
$ws.Range("E2").Value = "'  +0.13%  "
$ws.Range("D3").Value = "'3.372.37"
$ws.Range("E3").Value = "'  -3.21%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'580.48"
$ws.Range("E5").Value = "'  -0.38%  "
$ws.Range("D6").Value = "'180.03"
$ws.Range("E6").Value = "'  +3.51%  "
$ws.Range("E7").Value = "'  +4.95%  "
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("D9").Value = "'3.370.06"
$ws.Range("E9").Value = "'  -3.35%  "
$ws.Range("E10").Value = "'  -0.65%  "
$ws.Range("D11").Value = "'6.93"
$ws.Range("E11").Value = "'  +1.10%  "
$ws.Range("E12").Value = "'  +0.58%  "
$ws.Range("D13").Value = "'3.959.43"
$ws.Range("E13").Value = "'  -3.11%  "
$ws.Range("D15").Value = "'28.94"
$ws.Range("E15").Value = "'  -3.82%  "
$ws.Range("D16").Value = "'65.974.01"
$ws.Range("E16").Value = "'  -0.15%  "
$ws.Range("E17").Value = "'  -0.37%  "
$ws.Range("D18").Value = "'3.373.73"
$ws.Range("E18").Value = "'  -3.05%  "
$ws.Range("E19").Value = "'  -2.29%  "
$ws.Range("D20").Value = "'13.60"
$ws.Range("E20").Value = "'  -2.34%  "
$ws.Range("D21").Value = "'365.73"
$ws.Range("E21").Value = "'  +0.01%  "
$ws.Range("E22").Value = "'  -3.11%  "
$ws.Range("D23").Value = "'72.36"
$ws.Range("E23").Value = "'  -0.51%  "
$ws.Range("E24").Value = "'  -0.32%  "
$ws.Range("E25").Value = "'  -1.26%  "
$ws.Range("E26").Value = "'  -0.21%  "
$ws.Range("E27").Value = "'  +0.58%  "
$ws.Range("E28").Value = "'  +0.50%  "
$ws.Range("E29").Value = "'  -0.06%  "
$ws.Range("B30").Value = "'NEARProtocol"
$ws.Range("C30").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'5.74"
$ws.Range("E30").Value = "'  -0.75%  "
$ws.Range("B31").Value = "'PancakeSwap"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.98"
$ws.Range("E31").Value = "'  -0.38%  "
$ws.Range("E32").Value = "'  -4.44%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "'  +0.00%  "
$ws.Range("E34").Value = "'  -2.29%  "
$ws.Range("E35").Value = "'  -3.88%  "
$ws.Range("D36").Value = "'1.52"
$ws.Range("E36").Value = "'  -1.59%  "
$ws.Range("D37").Value = "'160.94"
$ws.Range("E37").Value = "'  +0.67%  "
$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "'  -3.85%  "
$ws.Range("D39").Value = "'27.12"
$ws.Range("E39").Value = "'  -7.55%  "
$ws.Range("E40").Value = "'  +0.17%  "
$ws.Range("E41").Value = "'  +1.25%  "
$ws.Range("B42").Value = "'RenderToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'6.32"
$ws.Range("E42").Value = "'  -1.11%  "
$ws.Range("B43").Value = "'Maker"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.681.54"
$ws.Range("E43").Value = "'  -4.91%  "
$ws.Range("E44").Value = "'  -2.85%  "
$ws.Range("D45").Value = "'0.0675"
$ws.Range("E45").Value = "'  -1.58%  "
$ws.Range("D46").Value = "'338.66"
$ws.Range("E46").Value = "'  +10.37%  "
$ws.Range("D47").Value = "'39.91"
$ws.Range("E47").Value = "'  -0.17%  "
$ws.Range("D48").Value = "'24.37"
$ws.Range("E48").Value = "'  +0.86%  "
$ws.Range("D49").Value = "'0.0283"
$ws.Range("E49").Value = "'  -1.75%  "
$ws.Range("E50").Value = "'  +3.17%  "
$ws.Range("E51").Value = "'  +0.52%  "
